$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# --- Update "Carga de Catalogo PROBING" table (rows 3-5: Factor de Carga (PROBING)) ---
$ws.Range("B3").Value = 615801.42500000005
$ws.Range("C3").Value = 17700.328000000001
$ws.Range("B4").Value = 615801.42500000005
$ws.Range("C4").Value = 17489.276000000002
$ws.Range("B5").Value = 615801.42500000005
$ws.Range("C5").Value = 18044.236000000001

# --- Update "Carga de Catalogo CHAINING" table (rows 10-12: Factor de Carga (CHAINING)) ---
$ws.Range("B10").Value = 615821.24300000002
$ws.Range("C10").Value = 18932.419999999998
$ws.Range("B11").Value = 615821.24300000002
$ws.Range("C11").Value = 19604.988000000001
$ws.Range("B12").Value = 615821.18000000005
$ws.Range("C12").Value = 19057.512999999999

# --- Update the selected cell to match the author's final cursor position ---
$ws.Range("B13").Select()
